# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-03 07:15:12
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# --- helper: write a literal text value into a cell without Excel's
# "looks like a number/percentage" auto-conversion kicking in. We stage the
# text in a scratch cell that is explicitly Text-formatted, copy it, and
# paste only the *value* into the destination - this leaves the
# destination's own number format/style completely untouched. ---
function Set-LiteralText($range, [string]$text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

# G2: reorder the "Recorded By" email list (same members, different order)
Set-LiteralText $ws.Range("G2") "Amira.Sobhy@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg"

# Summary counters (K/L columns)
$ws.Range("L6").Value = 4          # Recorded Sessions
$ws.Range("L8").Value = 25         # Pending Sessions
Set-LiteralText $ws.Range("L9") "13.8%"    # Coverage %
Set-LiteralText $ws.Range("L10") "23.8%"   # Average Attendance %

# Matching per-class statistics block (row 15, columns O/Q/R/S)
$ws.Range("O15").Value = 4
$ws.Range("Q15").Value = 25
Set-LiteralText $ws.Range("R15") "13.8%"
Set-LiteralText $ws.Range("S15") "23.8%"

# Row 26 (Pharmacology C1 session 1) moved from Pending -> Recorded:
# pick up the "Recorded" look by copying the format from an existing
# Recorded row (row 9), then fill in the newly-recorded attendance info.
$ws.Range("A9:I9").Copy()
$ws.Range("A26:I26").PasteSpecial(-4122)   # xlPasteFormats

Set-LiteralText $ws.Range("G26") "nancy.abdelshafy@med.asu.edu.eg"
Set-LiteralText $ws.Range("H26") "102/251"
Set-LiteralText $ws.Range("I26") "Recorded"

$excel.CutCopyMode = 0
